$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "12-month period ended" year headers: drop oldest (1396/12), shift others left, append newest (1401/12)
$yearE = "دوازده ماهه منتهی به 1397/12"
$yearF = "دوازده ماهه منتهی به 1398/12"
$yearG = "دوازده ماهه منتهی به 1399/12"
$yearH = "دوازده ماهه منتهی به 1400/12"
$yearI = "دوازده ماهه منتهی به 1401/12"
$headerRows = @(8,27,38,49,60,71,82,93,104,115,124,134,144,153)
foreach ($r in $headerRows) {
  $ws.Range("E$r").Value = $yearE
  $ws.Range("F$r").Value = $yearF
  $ws.Range("G$r").Value = $yearG
  $ws.Range("H$r").Value = $yearH
  $ws.Range("I$r").Value = $yearI
}

# Shift each data row left by one year column (drop oldest year, shift others, append newest year value)
$ws.Range("E10").Value = 457086 ; $ws.Range("F10").Value = 539583 ; $ws.Range("G10").Value = 1621403 ; $ws.Range("H10").Value = 4536276 ; $ws.Range("I10").Value = 7254705
$ws.Range("E11").Value = 37591 ; $ws.Range("F11").Value = 45014 ; $ws.Range("G11").Value = 82363 ; $ws.Range("H11").Value = 142448 ; $ws.Range("I11").Value = 204289
$ws.Range("E12").Value = 281857 ; $ws.Range("F12").Value = 369393 ; $ws.Range("G12").Value = 538641 ; $ws.Range("H12").Value = 1084669 ; $ws.Range("I12").Value = 1669182
$ws.Range("E13").Value = 776534 ; $ws.Range("F13").Value = 953990 ; $ws.Range("G13").Value = 2242407 ; $ws.Range("H13").Value = 5763393 ; $ws.Range("I13").Value = 9128176
$ws.Range("E14").Value = 0 ; $ws.Range("F14").Value = 0 ; $ws.Range("G14").Value = 0 ; $ws.Range("H14").Value = -26414 ; $ws.Range("I14").Value = -63776
$ws.Range("E15").Value = 776534 ; $ws.Range("F15").Value = 953990 ; $ws.Range("G15").Value = 2242407 ; $ws.Range("H15").Value = 5736979 ; $ws.Range("I15").Value = 9064400
$ws.Range("E16").Value = 0 ; $ws.Range("F16").Value = 0 ; $ws.Range("G16").Value = 0 ; $ws.Range("H16").Value = 0 ; $ws.Range("I16").Value = 0
$ws.Range("E17").Value = 0 ; $ws.Range("F17").Value = 0 ; $ws.Range("G17").Value = 0 ; $ws.Range("H17").Value = 0 ; $ws.Range("I17").Value = 0
$ws.Range("E18").Value = 776534 ; $ws.Range("F18").Value = 953990 ; $ws.Range("G18").Value = 2242407 ; $ws.Range("H18").Value = 5736979 ; $ws.Range("I18").Value = 9064400
$ws.Range("E19").Value = 3102 ; $ws.Range("F19").Value = 23123 ; $ws.Range("G19").Value = 81114 ; $ws.Range("H19").Value = 567216 ; $ws.Range("I19").Value = 1682113
$ws.Range("E20").Value = -23123 ; $ws.Range("F20").Value = -81114 ; $ws.Range("G20").Value = -567216 ; $ws.Range("H20").Value = -1682113 ; $ws.Range("I20").Value = -3564885
$ws.Range("E21").Value = 756513 ; $ws.Range("F21").Value = 895999 ; $ws.Range("G21").Value = 1756305 ; $ws.Range("H21").Value = 4622082 ; $ws.Range("I21").Value = 7181628
$ws.Range("E22").Value = 0 ; $ws.Range("F22").Value = 0 ; $ws.Range("G22").Value = 2938 ; $ws.Range("H22").Value = 0 ; $ws.Range("I22").Value = 0
$ws.Range("E23").Value = 756513 ; $ws.Range("F23").Value = 895999 ; $ws.Range("G23").Value = 1759243 ; $ws.Range("H23").Value = 4622082 ; $ws.Range("I23").Value = 7181628
$ws.Range("E29").Value = "-" ; $ws.Range("F29").Value = "-" ; $ws.Range("G29").Value = 0 ; $ws.Range("H29").Value = 8040 ; $ws.Range("I29").Value = 160
$ws.Range("E30").Value = 0 ; $ws.Range("F30").Value = 0 ; $ws.Range("G30").Value = 0 ; $ws.Range("H30").Value = 0 ; $ws.Range("I30").Value = 0
$ws.Range("E31").Value = 401 ; $ws.Range("F31").Value = 584 ; $ws.Range("G31").Value = 541 ; $ws.Range("H31").Value = 396 ; $ws.Range("I31").Value = 637
$ws.Range("E32").Value = 127 ; $ws.Range("F32").Value = 125 ; $ws.Range("G32").Value = 138 ; $ws.Range("H32").Value = 42 ; $ws.Range("I32").Value = 90
$ws.Range("E33").Value = 852 ; $ws.Range("F33").Value = 520 ; $ws.Range("G33").Value = 1317 ; $ws.Range("H33").Value = 194 ; $ws.Range("I33").Value = 889
$ws.Range("E34").Value = 1380 ; $ws.Range("F34").Value = 1229 ; $ws.Range("G34").Value = 1996 ; $ws.Range("H34").Value = 8672 ; $ws.Range("I34").Value = 1776
$ws.Range("E40").Value = "-" ; $ws.Range("F40").Value = "-" ; $ws.Range("G40").Value = 14972 ; $ws.Range("H40").Value = 6136 ; $ws.Range("I40").Value = 15035
$ws.Range("E41").Value = 4186 ; $ws.Range("F41").Value = 9634 ; $ws.Range("G41").Value = 8753 ; $ws.Range("H41").Value = 7859 ; $ws.Range("I41").Value = 14205
$ws.Range("E42").Value = 18244 ; $ws.Range("F42").Value = 14510 ; $ws.Range("G42").Value = 14909 ; $ws.Range("H42").Value = 15949 ; $ws.Range("I42").Value = 11611
$ws.Range("E43").Value = 8344 ; $ws.Range("F43").Value = 8966 ; $ws.Range("G43").Value = 9054 ; $ws.Range("H43").Value = 8096 ; $ws.Range("I43").Value = 7646
$ws.Range("E44").Value = 15939 ; $ws.Range("F44").Value = 15105 ; $ws.Range("G44").Value = 14963 ; $ws.Range("H44").Value = 15441 ; $ws.Range("I44").Value = 15484
$ws.Range("E45").Value = 46713 ; $ws.Range("F45").Value = 48215 ; $ws.Range("G45").Value = 62651 ; $ws.Range("H45").Value = 53481 ; $ws.Range("I45").Value = 63981
$ws.Range("E51").Value = "-" ; $ws.Range("F51").Value = "-" ; $ws.Range("G51").Value = 6931 ; $ws.Range("H51").Value = 14016 ; $ws.Range("I51").Value = 14864
$ws.Range("E52").Value = 4186 ; $ws.Range("F52").Value = 9634 ; $ws.Range("G52").Value = 8753 ; $ws.Range("H52").Value = 7859 ; $ws.Range("I52").Value = 14205
$ws.Range("E53").Value = 18060 ; $ws.Range("F53").Value = 14553 ; $ws.Range("G53").Value = 15054 ; $ws.Range("H53").Value = 15708 ; $ws.Range("I53").Value = 11383
$ws.Range("E54").Value = 8346 ; $ws.Range("F54").Value = 8953 ; $ws.Range("G54").Value = 9150 ; $ws.Range("H54").Value = 8048 ; $ws.Range("I54").Value = 7619
$ws.Range("E55").Value = 16200 ; $ws.Range("F55").Value = 14308 ; $ws.Range("G55").Value = 16087 ; $ws.Range("H55").Value = 14746 ; $ws.Range("I55").Value = 13934
$ws.Range("E56").Value = 46792 ; $ws.Range("F56").Value = 47448 ; $ws.Range("G56").Value = 55975 ; $ws.Range("H56").Value = 60377 ; $ws.Range("I56").Value = 62005
$ws.Range("E62").Value = "-" ; $ws.Range("F62").Value = "-" ; $ws.Range("G62").Value = 8041 ; $ws.Range("H62").Value = 160 ; $ws.Range("I62").Value = 331
$ws.Range("E63").Value = 0 ; $ws.Range("F63").Value = 0 ; $ws.Range("G63").Value = 0 ; $ws.Range("H63").Value = 0 ; $ws.Range("I63").Value = 0
$ws.Range("E64").Value = 585 ; $ws.Range("F64").Value = 541 ; $ws.Range("G64").Value = 396 ; $ws.Range("H64").Value = 637 ; $ws.Range("I64").Value = 865
$ws.Range("E65").Value = 125 ; $ws.Range("F65").Value = 138 ; $ws.Range("G65").Value = 42 ; $ws.Range("H65").Value = 90 ; $ws.Range("I65").Value = 117
$ws.Range("E66").Value = 591 ; $ws.Range("F66").Value = 1317 ; $ws.Range("G66").Value = 193 ; $ws.Range("H66").Value = 889 ; $ws.Range("I66").Value = 2439
$ws.Range("E67").Value = 1301 ; $ws.Range("F67").Value = 1996 ; $ws.Range("G67").Value = 8672 ; $ws.Range("H67").Value = 1776 ; $ws.Range("I67").Value = 3752
$ws.Range("E73").Value = "-" ; $ws.Range("F73").Value = "-" ; $ws.Range("G73").Value = 0 ; $ws.Range("H73").Value = 893337 ; $ws.Range("I73").Value = 17840
$ws.Range("E74").Value = 0 ; $ws.Range("F74").Value = 0 ; $ws.Range("G74").Value = 0 ; $ws.Range("H74").Value = 0 ; $ws.Range("I74").Value = 0
$ws.Range("E75").Value = 1793 ; $ws.Range("F75").Value = 4137 ; $ws.Range("G75").Value = 1703 ; $ws.Range("H75").Value = 9645 ; $ws.Range("I75").Value = 42072
$ws.Range("E76").Value = 1596 ; $ws.Range("F76").Value = 1715 ; $ws.Range("G76").Value = 2883 ; $ws.Range("H76").Value = 1566 ; $ws.Range("I76").Value = 11885
$ws.Range("E77").Value = 10358 ; $ws.Range("F77").Value = 9877 ; $ws.Range("G77").Value = 41534 ; $ws.Range("H77").Value = 10296 ; $ws.Range("I77").Value = 175750
$ws.Range("E78").Value = 13747 ; $ws.Range("F78").Value = 15729 ; $ws.Range("G78").Value = 46120 ; $ws.Range("H78").Value = 914844 ; $ws.Range("I78").Value = 247547
$ws.Range("E84").Value = "-" ; $ws.Range("F84").Value = "-" ; $ws.Range("G84").Value = 1493688 ; $ws.Range("H84").Value = 739853 ; $ws.Range("I84").Value = 2964620
$ws.Range("E85").Value = 9002 ; $ws.Range("F85").Value = 19170 ; $ws.Range("G85").Value = 20719 ; $ws.Range("H85").Value = 66700 ; $ws.Range("I85").Value = 313219
$ws.Range("E86").Value = 111551 ; $ws.Range("F86").Value = 75968 ; $ws.Range("G86").Value = 150788 ; $ws.Range("H86").Value = 615039 ; $ws.Range("I86").Value = 420739
$ws.Range("E87").Value = 115773 ; $ws.Range("F87").Value = 134246 ; $ws.Range("G87").Value = 216100 ; $ws.Range("H87").Value = 663051 ; $ws.Range("I87").Value = 1498111
$ws.Range("E88").Value = 222741 ; $ws.Range("F88").Value = 340589 ; $ws.Range("G88").Value = 608832 ; $ws.Range("H88").Value = 1784337 ; $ws.Range("I88").Value = 2191749
$ws.Range("E89").Value = 459067 ; $ws.Range("F89").Value = 569973 ; $ws.Range("G89").Value = 2490127 ; $ws.Range("H89").Value = 3868980 ; $ws.Range("I89").Value = 7388438
$ws.Range("E95").Value = "-" ; $ws.Range("F95").Value = "-" ; $ws.Range("G95").Value = 600351 ; $ws.Range("H95").Value = 1615349 ; $ws.Range("I95").Value = 2900036
$ws.Range("E96").Value = 9002 ; $ws.Range("F96").Value = 19170 ; $ws.Range("G96").Value = 20719 ; $ws.Range("H96").Value = 66700 ; $ws.Range("I96").Value = 313219
$ws.Range("E97").Value = 109207 ; $ws.Range("F97").Value = 78402 ; $ws.Range("G97").Value = 142845 ; $ws.Range("H97").Value = 582613 ; $ws.Range("I97").Value = 434249
$ws.Range("E98").Value = 115656 ; $ws.Range("F98").Value = 133079 ; $ws.Range("G98").Value = 217417 ; $ws.Range("H98").Value = 652731 ; $ws.Range("I98").Value = 1491307
$ws.Range("E99").Value = 223221 ; $ws.Range("F99").Value = 308932 ; $ws.Range("G99").Value = 640071 ; $ws.Range("H99").Value = 1618883 ; $ws.Range("I99").Value = 2115894
$ws.Range("E100").Value = 457086 ; $ws.Range("F100").Value = 539583 ; $ws.Range("G100").Value = 1621403 ; $ws.Range("H100").Value = 4536276 ; $ws.Range("I100").Value = 7254705
$ws.Range("E106").Value = "-" ; $ws.Range("F106").Value = "-" ; $ws.Range("G106").Value = 893337 ; $ws.Range("H106").Value = 17841 ; $ws.Range("I106").Value = 82424
$ws.Range("E107").Value = 0 ; $ws.Range("F107").Value = 0 ; $ws.Range("G107").Value = 0 ; $ws.Range("H107").Value = 0 ; $ws.Range("I107").Value = 0
$ws.Range("E108").Value = 4137 ; $ws.Range("F108").Value = 1703 ; $ws.Range("G108").Value = 9646 ; $ws.Range("H108").Value = 42071 ; $ws.Range("I108").Value = 28562
$ws.Range("E109").Value = 1713 ; $ws.Range("F109").Value = 2882 ; $ws.Range("G109").Value = 1566 ; $ws.Range("H109").Value = 11886 ; $ws.Range("I109").Value = 18689
$ws.Range("E110").Value = 9878 ; $ws.Range("F110").Value = 41534 ; $ws.Range("G110").Value = 10295 ; $ws.Range("H110").Value = 175750 ; $ws.Range("I110").Value = 251605
$ws.Range("E111").Value = 15728 ; $ws.Range("F111").Value = 46119 ; $ws.Range("G111").Value = 914844 ; $ws.Range("H111").Value = 247548 ; $ws.Range("I111").Value = 381280
$ws.Range("E117").Value = "-" ; $ws.Range("F117").Value = "-" ; $ws.Range("G117").Value = "-" ; $ws.Range("H117").Value = 111111567 ; $ws.Range("I117").Value = 111500000
$ws.Range("E118").Value = 4471322 ; $ws.Range("F118").Value = 7083904 ; $ws.Range("G118").Value = 3147874 ; $ws.Range("H118").Value = 24356061 ; $ws.Range("I118").Value = 66047096
$ws.Range("E119").Value = 12566929 ; $ws.Range("F119").Value = 13720000 ; $ws.Range("G119").Value = 20891304 ; $ws.Range("H119").Value = 37285714 ; $ws.Range("I119").Value = 132055556
$ws.Range("E120").Value = 12157277 ; $ws.Range("F120").Value = 18994231 ; $ws.Range("G120").Value = 31536826 ; $ws.Range("H120").Value = 53072165 ; $ws.Range("I120").Value = 197694038
$ws.Range("E126").Value = "-" ; $ws.Range("F126").Value = "-" ; $ws.Range("G126").Value = 99765429 ; $ws.Range("H126").Value = 120575782 ; $ws.Range("I126").Value = 197181244
$ws.Range("E127").Value = 2150502 ; $ws.Range("F127").Value = 1989828 ; $ws.Range("G127").Value = 2367074 ; $ws.Range("H127").Value = 8487085 ; $ws.Range("I127").Value = 22049912
$ws.Range("E128").Value = 6114394 ; $ws.Range("F128").Value = 5235562 ; $ws.Range("G128").Value = 10113891 ; $ws.Range("H128").Value = 38562857 ; $ws.Range("I128").Value = 36236242
$ws.Range("E129").Value = 13875000 ; $ws.Range("F129").Value = 14972786 ; $ws.Range("G129").Value = 23867904 ; $ws.Range("H129").Value = 81898592 ; $ws.Range("I129").Value = 195933952
$ws.Range("E130").Value = 13974591 ; $ws.Range("F130").Value = 22548097 ; $ws.Range("G130").Value = 40689167 ; $ws.Range("H130").Value = 115558384 ; $ws.Range("I130").Value = 141549277
$ws.Range("E136").Value = "-" ; $ws.Range("F136").Value = "-" ; $ws.Range("G136").Value = 86618237 ; $ws.Range("H136").Value = 115250357 ; $ws.Range("I136").Value = 195104682
$ws.Range("E137").Value = 2150502 ; $ws.Range("F137").Value = 1989828 ; $ws.Range("G137").Value = 2367074 ; $ws.Range("H137").Value = 8487085 ; $ws.Range("I137").Value = 22049912
$ws.Range("E138").Value = 6046899 ; $ws.Range("F138").Value = 5387343 ; $ws.Range("G138").Value = 9488840 ; $ws.Range("H138").Value = 37090209 ; $ws.Range("I138").Value = 38148906
$ws.Range("E139").Value = 13857656 ; $ws.Range("F139").Value = 14864180 ; $ws.Range("G139").Value = 23761421 ; $ws.Range("H139").Value = 81104747 ; $ws.Range("I139").Value = 195735267
$ws.Range("E140").Value = 13779074 ; $ws.Range("F140").Value = 21591557 ; $ws.Range("G140").Value = 39788090 ; $ws.Range("H140").Value = 109784552 ; $ws.Range("I140").Value = 151851155
$ws.Range("E146").Value = "-" ; $ws.Range("F146").Value = "-" ; $ws.Range("G146").Value = 111097749 ; $ws.Range("H146").Value = 111506250 ; $ws.Range("I146").Value = 249015106
$ws.Range("E147").Value = 7071795 ; $ws.Range("F147").Value = 3147874 ; $ws.Range("G147").Value = 24358586 ; $ws.Range("H147").Value = 66045526 ; $ws.Range("I147").Value = 33019653
$ws.Range("E148").Value = 13704000 ; $ws.Range("F148").Value = 20884058 ; $ws.Range("G148").Value = 37285714 ; $ws.Range("H148").Value = 132066667 ; $ws.Range("I148").Value = 159735043
$ws.Range("E149").Value = 16714044 ; $ws.Range("F149").Value = 31536826 ; $ws.Range("G149").Value = 53341969 ; $ws.Range("H149").Value = 197694038 ; $ws.Range("I149").Value = 103159082
$ws.Range("E159").Value = 5495 ; $ws.Range("F159").Value = 72295 ; $ws.Range("G159").Value = 132990 ; $ws.Range("H159").Value = 242458 ; $ws.Range("I159").Value = 216183
$ws.Range("E160").Value = 20502 ; $ws.Range("F160").Value = 26559 ; $ws.Range("G160").Value = 34696 ; $ws.Range("H160").Value = 163095 ; $ws.Range("I160").Value = 179820
$ws.Range("E161").Value = 57003 ; $ws.Range("F161").Value = 43649 ; $ws.Range("G161").Value = 52392 ; $ws.Range("H161").Value = 67075 ; $ws.Range("I161").Value = 61944
$ws.Range("E162").Value = 90216 ; $ws.Range("F162").Value = 144786 ; $ws.Range("G162").Value = 188024 ; $ws.Range("H162").Value = 268629 ; $ws.Range("I162").Value = 313267
$ws.Range("E163").Value = 0 ; $ws.Range("F163").Value = 0 ; $ws.Range("G163").Value = 0 ; $ws.Range("H163").Value = 0 ; $ws.Range("I163").Value = 0
$ws.Range("E164").Value = 108641 ; $ws.Range("F164").Value = 82104 ; $ws.Range("G164").Value = 130539 ; $ws.Range("H164").Value = 343412 ; $ws.Range("I164").Value = 897968
$ws.Range("E165").Value = 281857 ; $ws.Range("F165").Value = 369393 ; $ws.Range("G165").Value = 538641 ; $ws.Range("H165").Value = 1084669 ; $ws.Range("I165").Value = 1669182
